# Generate Report for Handback
# Updates the af1da157-1533-4987-85f2-0cc18e314ccf row (row 7) on the
# zh-cn and de-de sheets: the handback for that file has now come in, so
# the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns get filled in, and an "Error Detail" is recorded
# because the handback version is behind the latest source.

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/17a27556acf95d19a289e84dc2b6790989b3fb08/e2e/af1da157-1533-4987-85f2-0cc18e314ccf.md"
$targetDisplay = "af1da157-1533-4987-85f2-0cc18e314ccf.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/14c6b814c3e3647b8e03c74c683e4db6ec12cc4e/e2e/af1da157-1533-4987-85f2-0cc18e314ccf.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/17a27556acf95d19a289e84dc2b6790989b3fb08/e2e/af1da157-1533-4987-85f2-0cc18e314ccf.md."

# --- zh-cn sheet, row 7 ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $targetUrl, "", "", $targetDisplay)
$wsZh.Range("J7").Value = "af1da157-1533-4987-85f2-0cc18e314ccf.5c697aeb0a13452ccb5786a1ddff926bb3f9b14b.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-13 15:07:52"
$wsZh.Range("P7").Value = $errorDetail

# --- de-de sheet, row 7 ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $targetUrl, "", "", $targetDisplay)
$wsDe.Range("J7").Value = "af1da157-1533-4987-85f2-0cc18e314ccf.5c697aeb0a13452ccb5786a1ddff926bb3f9b14b.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-13 15:08:03"
$wsDe.Range("P7").Value = $errorDetail
